# Ajout de la fonctionnalite GANTT
# - hide the helper columns (E:G) used to compute the Gantt chart and
#   collapse their width to 0, hide the old "+TARD" column (H) too
#   (kept for calculations but no longer meant to be seen),
# - move the active selection to D5,
# - drop the one-off style that had been applied to J3 so it falls back
#   to the sheet's default formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TACHES")

# Hide columns E:G completely (width 0) - these hold the GANTT helper data.
$helperCols = $ws.Range("E1:G1").EntireColumn
$helperCols.ColumnWidth = -0.8333333333333334
$helperCols.Hidden = $true

# Column H keeps its width (14.5) but is now hidden as well.
$hCol = $ws.Range("H1").EntireColumn
$hCol.ColumnWidth = 13.666666666666666
$hCol.Hidden = $true

# J3 no longer needs its dedicated style - reset it to the default.
$ws.Range("J3").Style = "Normal"

# Move the selection to D5.
$ws.Range("D5").Select() | Out-Null
